$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row at row 3, shifting existing rows 3:74 down to 4:75
$ws.Rows.Item(3).Insert()

# Fill in the new row 3 with the new observation's data.
# Columns A, B, C, E-J are identical across the whole sheet (constant
# market/product metadata), so copy them down from row 2.
$ws.Range("A3").Value2 = $ws.Range("A2").Value2
$ws.Range("B3").Value2 = $ws.Range("B2").Value2
$ws.Range("C3").Value2 = $ws.Range("C2").Value2
$ws.Range("D3").Value2 = 44922
$ws.Range("E3").Value2 = $ws.Range("E2").Value2
$ws.Range("F3").Value2 = $ws.Range("F2").Value2
$ws.Range("G3").Value2 = $ws.Range("G2").Value2
$ws.Range("H3").Value2 = $ws.Range("H2").Value2
$ws.Range("I3").Value2 = $ws.Range("I2").Value2
$ws.Range("J3").Value2 = $ws.Range("J2").Value2
$ws.Range("K3").Value2 = "Black Amber"
$ws.Range("L3").Value2 = "Primera"
$ws.Range("M3").Value2 = 100
$ws.Range("N3").Value2 = 13000
$ws.Range("O3").Value2 = 14000
$ws.Range("P3").Value2 = 13500
$ws.Range("Q3").Value2 = "$/bandeja 18 kilos granel"
$ws.Range("R3").Value2 = "Región de O'Higgins"
$ws.Range("S3").Value2 = 750
$ws.Range("T3").Value2 = 18
